{"js": "// The paragraph \"<id>p045r_1</id>\" is currently split across three runs\n// (the run boundaries come from the tool that generated the original OOXML,\n// not from any formatting difference that should remain visible):\n//   1. \"<id>\"      - Courier New, color 7f6000, sz 18\n//   2. \"p045r_1\"   - color 000000\n//   3. \"</id>\"     - Courier New, color 7f6000, sz 18\n// The edit collapses these three runs into a single run containing the full\n// \"<id>p045r_1</id>\" text, carrying the formatting of the first run.\n\nconst searchText = \"<id>p045r_1</id>\";\n\nconst results = context.document.body.search(searchText, { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Could not find the target text '\" + searchText + \"' to merge.\");\n}\n\n// There is only one real match (the sibling \"<id>fig_p045r_1</id>\" paragraph\n// has different text), but guard against unexpected duplicates by only\n// touching the first hit.\nconst target = results.items[0];\n\n// Re-inserting the same text over the matched range merges the runs that\n// span it into a single run, taking on the formatting of the range's start\n// (i.e. the first run's Courier New / 7f6000 / sz18 properties), which is\n// exactly the collapse described by the diff.\ntarget.insertText(searchText, \"Replace\");\n\nawait context.sync();\n", "ps1": "# The paragraph \"<id>p045r_1</id>\" is currently split across three runs\n# (the run boundaries come from the tool that generated the original OOXML,\n# not from any formatting difference that should remain visible):\n#   1. \"<id>\"      - Courier New, color 7f6000, sz 18\n#   2. \"p045r_1\"   - color 000000\n#   3. \"</id>\"     - Courier New, color 7f6000, sz 18\n# The edit collapses these three runs into a single run containing the full\n# \"<id>p045r_1</id>\" text, carrying the formatting of the first run.\n\n$d = $word.ActiveDocument\n\n$fullText = \"<id>p045r_1</id>\"\n$firstRunText = \"<id>\"\n\n# Locate the (unique) span containing the full text - a plain, non-wildcard\n# search so \"<\", \">\", \"/\" are matched literally. This only moves the Find\n# range; passing \"\" / wdReplaceNone(0) performs no replacement/mutation yet.\n$rng = $d.Content\n$rng.Find.ClearFormatting()\n$found = $rng.Find.Execute($fullText, $false, $false, $false, $false, $false, $true, 1, $false, \"\", 0)\nif (-not $found) {\n    throw \"Could not find the target text '$fullText' to merge.\"\n}\n\n$matchStart = $rng.Start\n$matchEnd = $rng.End\n$splitPoint = $matchStart + $firstRunText.Length\n\n# Delete the 2nd/3rd runs' content (\"p045r_1</id>\") out of the match, then\n# append it onto the end of the surviving first run (\"<id>\"). Growing the\n# first run's own range (instead of rewriting the whole paragraph) keeps\n# that run's identity/formatting - Courier New / 7f6000 / sz18 - and lets\n# the trailing empty run stay untouched, exactly mirroring the diff's merge\n# of the three runs into one.\n$remainderRange = $d.Range($splitPoint, $matchEnd)\n$remainderText = $remainderRange.Text\n$remainderRange.Delete()\n\n$firstRun = $d.Range($matchStart, $splitPoint)\n$firstRun.InsertAfter($remainderText)\n"}
